$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.989.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.925.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.60%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.42"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.35%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.53"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.73%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.508"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.93%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.91"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.81%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.144"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.30%  "

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.15%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.60"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.72%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.24%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.412.41"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.935.57"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.37%  "

# Row 17
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.71"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.36%  "

# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.925.49"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.51%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "432.31"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.39"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.48%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.679"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.38%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.09"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.58%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.51"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.48%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.00%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.20"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.34%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.88"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.07%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.17%  "

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.24%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.19%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.66"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.77%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.49%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.05%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0857"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.99%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.54%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.64"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.12%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.05%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.99"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.63%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.122"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.19%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.58"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.85%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.93%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.01"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.32%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "380.87"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.63%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.701.13"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.27%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0343"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.29%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.84"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.85%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.83"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.91%  "

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.29%  "

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.78%  "

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.62%  "
